# Commit: "completed o__bacillales, o__Erysipelotrichales, o__Lactobacillales,
#          o__mycoplasmatales, o__RFN20"
#
# For o__Lactobacillales.xlsx this adds a second worksheet
# "o__Lactobacillales_pred-t-p" (the "-p" = p-values/extra columns sheet)
# right after the existing "o__Lactobacillales_pred-t" sheet, carrying the
# same family columns plus max / prediction / rejection-f / gtdb-Tk columns.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after the existing one so it becomes the active tab,
# matching <workbookView ... activeTab="1"/> and sheet order in the target file.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o__Lactobacillales_pred-t-p"

# Match sheet 1's page margins (0.75/0.75/1/1/0.5/0.5 in) instead of the
# worksheet-add default (0.7/0.7/0.75/0.75/0.3/0.3 in).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# ---- Header row -----------------------------------------------------
$ws2.Range("B1").Value = "f__Aerococcaceae"
$ws2.Range("C1").Value = "f__Carnobacteriaceae"
$ws2.Range("D1").Value = "f__Enterococcaceae"
$ws2.Range("E1").Value = "f__Lactobacillaceae"
$ws2.Range("F1").Value = "f__Listeriaceae"
$ws2.Range("G1").Value = "f__Streptococcaceae"
$ws2.Range("H1").Value = "max"
$ws2.Range("I1").Value = "prediction"
$ws2.Range("J1").Value = "rejection-f"
$ws2.Range("K1").Value = "gtdb-Tk"

# ---- Row labels -------------------------------------------------------
$ws2.Range("A2").Value = "RUG513"
$ws2.Range("A3").Value = "RUG664"

# ---- RUG513 data row ----------------------------------------------------
$ws2.Range("B2").Value = 0.00001031677749664935
$ws2.Range("C2").Value = 0.0012076969071859859
$ws2.Range("D2").Value = 0.00069296531361251554
$ws2.Range("E2").Value = 0.019376573798247441
$ws2.Range("F2").Value = 0.0091390052940975384
$ws2.Range("G2").Value = 0.96957344190935979
$ws2.Range("H2").Value = 0.96957344190935979
$ws2.Range("I2").Value = "f__Streptococcaceae"
$ws2.Range("J2").Value = "f__Streptococcaceae"
$ws2.Range("K2").Value = "f__Streptococcaceae"

# ---- RUG664 data row ----------------------------------------------------
$ws2.Range("B3").Value = 0.0000007360054568452776
$ws2.Range("C3").Value = 0.00078725939660665377
$ws2.Range("D3").Value = 0.0007426847517356998
$ws2.Range("E3").Value = 0.01270553146548673
$ws2.Range("F3").Value = 0.0041139014905610559
$ws2.Range("G3").Value = 0.98164988689015309
$ws2.Range("H3").Value = 0.98164988689015309
$ws2.Range("I3").Value = "f__Streptococcaceae"
$ws2.Range("J3").Value = "f__Streptococcaceae"
$ws2.Range("K3").Value = "f__Streptococcaceae"

# ---- Formatting: reuse the bordered/bold/centered style already used on
# sheet 1's header row + row labels (cell B1 there carries it) instead of
# creating a brand-new style entry. ------------------------------------
$ws1.Range("B1").Copy() | Out-Null
$ws2.Range("B1:K1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false | Out-Null

$ws2.Range("A1").Select() | Out-Null
